# Implemented spatial rollback, with some supporting changes in other
# components: the "rollback" / "inventory_year" documentation row gains a
# new option (pointing the field at a full layer definition rather than
# requiring a plain path), so a new row must be inserted into the
# config-parsing reference table on Sheet1, pushing the existing
# "required field" / "rollback_year" rows down by one.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the current row 58 ("rollback"/"inventory_year"
# null/empty requirement), shifting rows 58-60 down to 59-61.
$ws.Rows(58).Insert()

# Populate the newly-inserted row 58 with the new documentation entry:
# inventory_year can now also be given as a full <layer definition>.
$ws.Range("A58").Value = "rollback"
$ws.Range("B58").Value = "inventory_year"
$ws.Range("C58").Value = "Y"
$ws.Range("D58").Value = "<layer definition>"
$ws.Range("E58").Value = "Use as full layer definition for the inventory vintage layer."
$ws.Range("F58").Value = "Same as layer definition validation."

# Match the author's final selection/scroll position in the sheet.
$ws.Range("A59").Select()
